# Wellness.xlsx update: append 10 new training-log rows (223-232, date 08/28/2025)
# and register a new "Cheville " localisation value used by one of them.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 222 already carries the correct per-column number formats/styles
# (date style on A, "Helvetica Neue" style on B-H, centered-empty style on G).
# Copying it down first means every new row inherits that formatting for free;
# the actual data (and the I-column formula) is then written on top of it.
$templateRow = 222

# G206 already holds free-text ("Mollet ") in the "Localisation douleur" column,
# i.e. it carries the non-centered, left-aligned text style used whenever that
# column is populated (as opposed to the s=2 style used when G is left blank).
$gTextStyleCell = $ws.Cells.Item(206, 7)

# --- Row 223: Sofiane Belle ---
$ws.Range("A" + $templateRow + ":I" + $templateRow).Copy($ws.Range("A223:I223"))
$ws.Cells.Item(223, 1).Value = 45897        # Date
$ws.Cells.Item(223, 2).Value = "Sofiane Belle"
$ws.Cells.Item(223, 3).Value = 45         # Volume
$ws.Cells.Item(223, 4).Value = 2         # Intensité
$ws.Cells.Item(223, 5).Value = 3         # Fatigue
$ws.Cells.Item(223, 6).Value = 0         # Douleur
$ws.Cells.Item(223, 8).Value = 2         # Plaisir
$ws.Cells.Item(223, 9).Formula = "=C223*D223"   # Charge

# --- Row 224: Amir Etien ---
$ws.Range("A" + $templateRow + ":I" + $templateRow).Copy($ws.Range("A224:I224"))
$ws.Cells.Item(224, 1).Value = 45897        # Date
$ws.Cells.Item(224, 2).Value = "Amir Etien"
$ws.Cells.Item(224, 3).Value = 45         # Volume
$ws.Cells.Item(224, 4).Value = 4         # Intensité
$ws.Cells.Item(224, 5).Value = 6         # Fatigue
$ws.Cells.Item(224, 6).Value = 0         # Douleur
$ws.Cells.Item(224, 8).Value = 7         # Plaisir
$ws.Cells.Item(224, 9).Formula = "=C224*D224"   # Charge

# --- Row 225: Amir Kherrab ---
$ws.Range("A" + $templateRow + ":I" + $templateRow).Copy($ws.Range("A225:I225"))
$ws.Cells.Item(225, 1).Value = 45897        # Date
$ws.Cells.Item(225, 2).Value = "Amir Kherrab"
$ws.Cells.Item(225, 3).Value = 45         # Volume
$ws.Cells.Item(225, 4).Value = 3         # Intensité
$ws.Cells.Item(225, 5).Value = 4         # Fatigue
$ws.Cells.Item(225, 6).Value = 0         # Douleur
$ws.Cells.Item(225, 8).Value = 6         # Plaisir
$ws.Cells.Item(225, 9).Formula = "=C225*D225"   # Charge

# --- Row 226: Emmanuel Valey ---
$ws.Range("A" + $templateRow + ":I" + $templateRow).Copy($ws.Range("A226:I226"))
$ws.Cells.Item(226, 1).Value = 45897        # Date
$ws.Cells.Item(226, 2).Value = "Emmanuel Valey"
$ws.Cells.Item(226, 3).Value = 45         # Volume
$ws.Cells.Item(226, 4).Value = 1         # Intensité
$ws.Cells.Item(226, 5).Value = 0         # Fatigue
$ws.Cells.Item(226, 6).Value = 3         # Douleur
$gTextStyleCell.Copy($ws.Cells.Item(226, 7))   # give G the "has text" style
$ws.Cells.Item(226, 7).Value = "Adducteur "   # Localisation douleur
$ws.Cells.Item(226, 8).Value = 2         # Plaisir
$ws.Cells.Item(226, 9).Formula = "=C226*D226"   # Charge

# --- Row 227: Karahali Souaré ---
$ws.Range("A" + $templateRow + ":I" + $templateRow).Copy($ws.Range("A227:I227"))
$ws.Cells.Item(227, 1).Value = 45897        # Date
$ws.Cells.Item(227, 2).Value = "Karahali Souaré"
$ws.Cells.Item(227, 3).Value = 45         # Volume
$ws.Cells.Item(227, 4).Value = 5         # Intensité
$ws.Cells.Item(227, 5).Value = 6         # Fatigue
$ws.Cells.Item(227, 6).Value = 8         # Douleur
$gTextStyleCell.Copy($ws.Cells.Item(227, 7))   # give G the "has text" style
$ws.Cells.Item(227, 7).Value = "Cheville "   # Localisation douleur
$ws.Cells.Item(227, 8).Value = 7         # Plaisir
$ws.Cells.Item(227, 9).Formula = "=C227*D227"   # Charge

# --- Row 228: Naim Dhib ---
$ws.Range("A" + $templateRow + ":I" + $templateRow).Copy($ws.Range("A228:I228"))
$ws.Cells.Item(228, 1).Value = 45897        # Date
$ws.Cells.Item(228, 2).Value = "Naim Dhib"
$ws.Cells.Item(228, 3).Value = 45         # Volume
$ws.Cells.Item(228, 4).Value = 3         # Intensité
$ws.Cells.Item(228, 5).Value = 3         # Fatigue
$ws.Cells.Item(228, 6).Value = 1         # Douleur
$gTextStyleCell.Copy($ws.Cells.Item(228, 7))   # give G the "has text" style
$ws.Cells.Item(228, 7).Value = "Courbature "   # Localisation douleur
$ws.Cells.Item(228, 8).Value = 6         # Plaisir
$ws.Cells.Item(228, 9).Formula = "=C228*D228"   # Charge

# --- Row 229: Yoan Zouma ---
$ws.Range("A" + $templateRow + ":I" + $templateRow).Copy($ws.Range("A229:I229"))
$ws.Cells.Item(229, 1).Value = 45897        # Date
$ws.Cells.Item(229, 2).Value = "Yoan Zouma"
$ws.Cells.Item(229, 3).Value = 45         # Volume
$ws.Cells.Item(229, 4).Value = 1         # Intensité
$ws.Cells.Item(229, 5).Value = 3         # Fatigue
$ws.Cells.Item(229, 6).Value = 0         # Douleur
$ws.Cells.Item(229, 8).Value = 1         # Plaisir
$ws.Cells.Item(229, 9).Formula = "=C229*D229"   # Charge

# --- Row 230: Yanis Berrached ---
$ws.Range("A" + $templateRow + ":I" + $templateRow).Copy($ws.Range("A230:I230"))
$ws.Cells.Item(230, 1).Value = 45897        # Date
$ws.Cells.Item(230, 2).Value = "Yanis Berrached"
$ws.Cells.Item(230, 3).Value = 45         # Volume
$ws.Cells.Item(230, 4).Value = 1         # Intensité
$ws.Cells.Item(230, 5).Value = 7         # Fatigue
$ws.Cells.Item(230, 6).Value = 0         # Douleur
$ws.Cells.Item(230, 8).Value = 0         # Plaisir
$ws.Cells.Item(230, 9).Formula = "=C230*D230"   # Charge

# --- Row 231: Naim Ighbane ---
$ws.Range("A" + $templateRow + ":I" + $templateRow).Copy($ws.Range("A231:I231"))
$ws.Cells.Item(231, 1).Value = 45897        # Date
$ws.Cells.Item(231, 2).Value = "Naim Ighbane"
$ws.Cells.Item(231, 3).Value = 45         # Volume
$ws.Cells.Item(231, 4).Value = 2         # Intensité
$ws.Cells.Item(231, 5).Value = 0         # Fatigue
$ws.Cells.Item(231, 6).Value = 0         # Douleur
$ws.Cells.Item(231, 8).Value = 5         # Plaisir
$ws.Cells.Item(231, 9).Formula = "=C231*D231"   # Charge

# --- Row 232: Ilan Ihaddadene ---
$ws.Range("A" + $templateRow + ":I" + $templateRow).Copy($ws.Range("A232:I232"))
$ws.Cells.Item(232, 1).Value = 45897        # Date
$ws.Cells.Item(232, 2).Value = "Ilan Ihaddadene"
$ws.Cells.Item(232, 3).Value = 45         # Volume
$ws.Cells.Item(232, 4).Value = 3         # Intensité
$ws.Cells.Item(232, 5).Value = 3         # Fatigue
$ws.Cells.Item(232, 6).Value = 0         # Douleur
$ws.Cells.Item(232, 8).Value = 7         # Plaisir
$ws.Cells.Item(232, 9).Formula = "=C232*D232"   # Charge

# Match the author's final scroll position/selection (sheet was scrolled down
# to keep showing the newly-entered rows, with K229 selected).
$win = $excel.ActiveWindow
$win.ScrollRow = 205
$win.ScrollColumn = 1
$ws.Range("K229").Select()

